$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row with Chinese labels (user can use Annotation to mark alias)
$ws.Range("A1").Value = "姓名"
$ws.Range("B1").Value = "昵称"
$ws.Range("C1").Value = "成绩"
$ws.Range("D1").Value = "年龄"

# Update selection to C2
$ws.Range("C2").Select()
